# Weekly update: insert a new latest-week record at row 3 (pushing the
# previously existing rows 3-18 down to 4-19) for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Cilantro".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (3..18) down by one to make room for the new entry.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's figures.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44532
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112040
$ws.Cells.Item(3, 7).Value = "Cilantro"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 2000
$ws.Cells.Item(3, 12).Value = 2200
$ws.Cells.Item(3, 13).Value = 2100
$ws.Cells.Item(3, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(3, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(3, 16).Value = 2100
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
